# Lab 5 - Completed upto Q4 Project - Handled services not in list
#
# The underlying data edit: 8 existing "Chat" rows in Sheet1 had their text
# corrected / replaced (typo fixes + a couple of rows repurposed into new,
# more generic fallback intents) so the chatbot gracefully handles
# "services not in list" queries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Order matters for how the shared-strings table gets rebuilt on save, so
# these are applied in the same sequence the original edit session used.
$ws.Range("A169").Value = "Available services please"
$ws.Range("A125").Value = "What are the service available on your saloon"
$ws.Range("A134").Value = "interjection usergreet I need to know about available saloon services"
$ws.Range("A166").Value = "interjection, Could you please provide info about services?"
$ws.Range("A223").Value = "Do your saloon provide services"
$ws.Range("A209").Value = "Do your saloon has "
$ws.Range("A130").Value = "Are there services"
$ws.Range("A136").Value = "Any services"

# Restore the view state (active selection) recorded in the saved workbook.
$ws.Activate() | Out-Null
$ws.Range("D212").Select() | Out-Null
